$d = $word.ActiveDocument

# The document's title banner references a docassemble variable with a
# single-underscore name, "interview_intro_prompt". The edit corrects it
# to the dotted attribute-access form "interview.intro_prompt" (mirroring
# "interview.title", used the same way later in the document).
$old = "interview_intro_prompt"
$new = "interview.intro_prompt"
$prefix = "interview."

# Locate the text to fix first (without replacing yet) so we can record
# where it starts.
$findRng = $d.Content
$findRng.Find.Execute($old, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$editStart = $findRng.Start

# Perform the actual correction.
$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                            $true, 1, $false, $new, 2)

if ($found) {
    # Word leaves its "_GoBack" bookmark (tracking the most recent edit
    # point) at the location of this change, right after the unchanged
    # "interview." prefix and before the corrected "intro_prompt" text.
    # Re-seating it there removes it from wherever it previously sat.
    $goBackStart = $editStart + $prefix.Length
    $goBackRange = $d.Range($goBackStart, $goBackStart)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
